$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '71.065.03'
$ws.Cells.Item(2, 5).Value = '  +0.26%  '

$ws.Cells.Item(3, 4).Value = '3.810.58'
$ws.Cells.Item(3, 5).Value = '  -0.86%  '

$ws.Cells.Item(4, 5).Value = '  -0.03%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '701.11'
$ws.Cells.Item(5, 5).Value = '  -0.45%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '172.25'
$ws.Cells.Item(6, 5).Value = '  -0.41%  '

$ws.Cells.Item(7, 4).Value = '3.809.34'
$ws.Cells.Item(7, 5).Value = '  -0.82%  '

$ws.Cells.Item(8, 5).Value = '  -0.04%  '

$ws.Cells.Item(9, 5).Value = '  +0.11%  '

$ws.Cells.Item(10, 5).Value = '  -0.56%  '

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '7.53'
$ws.Cells.Item(11, 5).Value = '  +2.59%  '

$ws.Cells.Item(12, 5).Value = '  +1.39%  '

$ws.Cells.Item(13, 5).Value = '  -1.18%  '

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '36.01'
$ws.Cells.Item(14, 5).Value = '  -1.34%  '

$ws.Cells.Item(15, 4).Value = '4.453.49'
$ws.Cells.Item(15, 5).Value = '  -0.85%  '

$ws.Cells.Item(16, 4).Value = '3.844.50'
$ws.Cells.Item(16, 5).Value = '  +0.23%  '

$ws.Cells.Item(17, 4).Value = '71.124.45'
$ws.Cells.Item(17, 5).Value = '  +0.26%  '

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '17.49'
$ws.Cells.Item(18, 5).Value = '  +0.96%  '

$ws.Cells.Item(19, 5).Value = '  -0.42%  '

$ws.Cells.Item(20, 5).Value = '  -0.54%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '512.99'
$ws.Cells.Item(21, 5).Value = '  +4.20%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '10.66'
$ws.Cells.Item(22, 5).Value = '  -0.16%  '

$ws.Cells.Item(23, 5).Value = '  -0.06%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '83.92'
$ws.Cells.Item(24, 5).Value = '  -1.47%  '

$ws.Cells.Item(26, 4).Value = '3.962.39'
$ws.Cells.Item(26, 5).Value = '  -0.93%  '

$ws.Cells.Item(27, 5).Value = '  -0.42%  '

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '10.43'
$ws.Cells.Item(28, 5).Value = '  -0.92%  '

$ws.Cells.Item(29, 5).Value = '  +0.03%  '

$ws.Cells.Item(30, 5).Value = '  -3.20%  '

$ws.Cells.Item(31, 5).Value = '  -4.22%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '7.44'
$ws.Cells.Item(32, 5).Value = '  -0.90%  '

$ws.Cells.Item(33, 5).Value = '  -1.46%  '

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '29.04'

$ws.Cells.Item(35, 5).Value = '  -4.61%  '

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '9.17'
$ws.Cells.Item(36, 5).Value = '  +0.26%  '

$ws.Cells.Item(37, 2).Value = 'Binance-PegBSC-USD'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '1.00'
$ws.Cells.Item(37, 5).Value = '  +0.10%  '

$ws.Cells.Item(38, 2).Value = 'RenzoRestakedETH'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Cells.Item(38, 4).Value = '3.772.33'
$ws.Cells.Item(38, 5).Value = '  -0.72%  '

$ws.Cells.Item(39, 5).Value = '  -2.08%  '

$ws.Cells.Item(40, 5).Value = '  +0.73%  '

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '6.02'
$ws.Cells.Item(41, 5).Value = '  -0.12%  '

$ws.Cells.Item(42, 5).Value = '  -0.71%  '

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '3.32'
$ws.Cells.Item(43, 5).Value = '  -0.03%  '

$ws.Cells.Item(44, 5).Value = '  -0.01%  '

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '173.14'
$ws.Cells.Item(45, 5).Value = '  +6.09%  '

$ws.Cells.Item(46, 5).Value = '  +0.09%  '

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '0.000315'
$ws.Cells.Item(47, 5).Value = '  +1.35%  '

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '431.32'
$ws.Cells.Item(48, 5).Value = '  +4.42%  '

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '49.39'
$ws.Cells.Item(49, 5).Value = '  +1.19%  '

$ws.Cells.Item(50, 5).Value = '  +0.07%  '

$ws.Cells.Item(51, 5).Value = '  +0.79%  '
